# "fix test 2 and update naming convention"
#
# The log had two "test"-phase entries on 3/30 (row 12) and 3/31 (row 13)
# that were logged against the wrong activity; the fix removes the
# duplicate/incorrect "Coding" entry on row 11 and the erroneous
# "Testing - Found some issues and the way to fix" entry on row 13,
# and re-times the remaining entries for that day (rows 7-10) onto 4/2/2019.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete bottom-most row first so the remaining row numbers below it
# don't shift out from under the next delete.
$ws.Rows("13:13").Delete()   # Testing / "Found some issues and the way to fix"
$ws.Rows("11:11").Delete()   # Coding / "Continue on coding" (duplicate entry)

# Rows 7-10 keep their activity/comment text, but get corrected dates and
# start/stop times.
$ws.Range("C7").Value = 43557
$ws.Range("D7").Value = 0.9194444444444444
$ws.Range("F7").Value = 0.92986111111111114

$ws.Range("C8").Value = 43557
$ws.Range("D8").Value = 0.92986111111111114
$ws.Range("F8").Value = 0.95972222222222225

$ws.Range("C9").Value = 43557
$ws.Range("D9").Value = 0.96250000000000002
$ws.Range("F9").Value = 0.98263888888888884

$ws.Range("C10").Value = 43557
$ws.Range("D10").Value = 0.98333333333333339
$ws.Range("F10").Value = 0.0041666666666666666

# Scroll/selection update left by the edit.
$ws.Range("C11").Select()
